# Fruta / hortaliza, semanal
#
# A new weekly price record is inserted as row 38 of the data table
# (pushing the former rows 38-92 down to 39-93). The new row carries a
# fresh observation for "Arándano (blue)" at Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 38..92 down to 39..93, opening up a blank row 38.
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new observation.
$ws.Range("A38").Value = 10
$ws.Range("B38").Value = "Vega Modelo de Temuco"
$ws.Range("C38").Value = "La Araucanía"
$ws.Range("D38").Value = 44671
$ws.Range("E38").Value = 9
$ws.Range("F38").Value = "Fruta"
$ws.Range("G38").Value = 100101
$ws.Range("H38").Value = "Berries"
$ws.Range("I38").Value = 100101001
$ws.Range("J38").Value = "Arándano (blue)"
$ws.Range("K38").Value = "Sin especificar"
$ws.Range("L38").Value = "Primera"
$ws.Range("M38").Value = 200
$ws.Range("N38").Value = 3000
$ws.Range("O38").Value = 3000
$ws.Range("P38").Value = 3000
$ws.Range("Q38").Value = "$/kilo"
$ws.Range("R38").Value = "Región de La Araucanía"
$ws.Range("S38").Value = 3000
$ws.Range("T38").Value = 1
